$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "306.34"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-3.58%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.39"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.71%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.053"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.93%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07583"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-6.78%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.250"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.81%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.586"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-9.50%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9043"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.95%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09972"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-11.11%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1765"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-4.89%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09005"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-3.36%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04389"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.18%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1054"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.15%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001231"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-6.32%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005836"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.42%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.369"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.51%"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.97%"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.85%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.846"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-7.61%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.22%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2851"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "9.77%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04160"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.23%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001219"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.94%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004069"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-4.76%"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.68%"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0002979"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "-0.22%"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02394"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-7.47%"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-6.60%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007857"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.31%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1301"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-6.72%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007117"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "9.29%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001955"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-6.50%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008384"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.83%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3322"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.44%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006435"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.49%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000753"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.13%"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-26.91%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.005773"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "70.12%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002107"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.13%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002007"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.13%"
